$wb = $excel.ActiveWorkbook

# --- Status text: "Ready for handoff" -> "In Translation" ---
# Overview sheet: per-language status columns (zh-cn -> E, de-de -> F) for both data rows.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# zh-cn sheet: Status column (C) for both data rows.
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "In Translation"
$zh.Range("C3").Value = "In Translation"

# de-de sheet: Status column (C) for both data rows.
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "In Translation"
$de.Range("C3").Value = "In Translation"

# --- Narrow the status columns ---
# Overview: columns E (zh-cn status) and F (de-de status).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de: column C (Status).
$zh.Columns.Item(3).ColumnWidth = 12.5
$de.Columns.Item(3).ColumnWidth = 12.5
